# Update the "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" sheet to reflect newly refreshed data (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (F column) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 457
$wsExpo.Range("F5").Value  = 1792
$wsExpo.Range("F6").Value  = 93
$wsExpo.Range("F7").Value  = 2255
$wsExpo.Range("F11").Value = 5094
$wsExpo.Range("F12").Value = 378
$wsExpo.Range("F17").Value = 208
$wsExpo.Range("F21").Value = 4136
$wsExpo.Range("F22").Value = 739
$wsExpo.Range("F23").Value = 744
$wsExpo.Range("F26").Value = 119
$wsExpo.Range("F30").Value = 103
$wsExpo.Range("F31").Value = 594
$wsExpo.Range("F33").Value = 30
$wsExpo.Range("F34").Value = 1052
$wsExpo.Range("F36").Value = 2667
$wsExpo.Range("F38").Value = 55

# --- Sheet "全部类型" (F column) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 457
$wsAll.Range("F5").Value  = 1792
$wsAll.Range("F6").Value  = 93
$wsAll.Range("F7").Value  = 2255
$wsAll.Range("F11").Value = 5094
$wsAll.Range("F12").Value = 378
$wsAll.Range("F17").Value = 208
$wsAll.Range("F21").Value = 4136
$wsAll.Range("F22").Value = 739
$wsAll.Range("F23").Value = 744
$wsAll.Range("F26").Value = 119
$wsAll.Range("F30").Value = 103
$wsAll.Range("F31").Value = 594
$wsAll.Range("F34").Value = 30
$wsAll.Range("F35").Value = 1052
$wsAll.Range("F37").Value = 2667
$wsAll.Range("F39").Value = 55
